$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("C2").Value = 0.5
$wsSummary.Range("D2").Value = 1
$wsSummary.Range("E2").Value = 0.6666666666666666
$wsSummary.Range("F2").Value = 0.8333333333333334
$wsSummary.Range("G2").Value = 0.9629629629629629
$wsSummary.Range("H2").Value = 0.5061650465008627
$wsSummary.Range("I2").Value = 534
$wsSummary.Range("J2").Value = 534
$wsSummary.Range("K2").Value = 0
$wsSummary.Range("L2").Value = 0

# --- Classification Report sheet ---
$wsClassification = $wb.Worksheets.Item("Classification Report")
$wsClassification.Range("B2").Value = 0
$wsClassification.Range("C2").Value = 0
$wsClassification.Range("D2").Value = 0
$wsClassification.Range("B3").Value = 0.5
$wsClassification.Range("C3").Value = 1
$wsClassification.Range("D3").Value = 0.6666666666666666

# --- Confusion Matrix sheet ---
$wsConfusion = $wb.Worksheets.Item("Confusion Matrix")
$wsConfusion.Range("B2").Value = 0
$wsConfusion.Range("C2").Value = 534
$wsConfusion.Range("B3").Value = 0
$wsConfusion.Range("C3").Value = 534
